# New tenant support in live
# Appends new sprint/interview-history rows to the AMSIN, BETA and AMS
# sheets (mirrors data typed into the live tracker) and normalises the
# formatting of the previous "last row" on AMSIN now that it is no longer
# the final entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Write literal text into a cell without Excel's "looks like a date"
# auto-conversion kicking in. We stage the literal text in a scratch
# cell (forced to Text format so it is never reinterpreted), copy it,
# and paste *values only* onto the destination - that leaves whatever
# number format/style the destination already carries untouched.
function Set-LiteralText($ws, $row, $col, $text) {
    $scratch = $ws.Cells.Item(600, 26)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
    $scratch.Clear()
}

# A plain (non date-like) text or numeric value - safe to assign
# directly, picks up the destination column's default styling.
function Set-PlainValue($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Clear()
    $ws.Cells.Item($row, $col).Value = $value
}

# The "Run Time" column: a date+time serial formatted with the sheet's
# custom "YYYY-MM-DD HH:MM:SS" number format.
function Set-RunTime($ws, $row, $serial) {
    $c = $ws.Cells.Item($row, 2)
    $c.Style = "Normal"
    $c.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $c.Value = $serial
}

# Writes one full data row: Run Date (text), Run Time (datetime),
# Sprint Name (text), Total/Pass/Fail cases, Time Taken.
function Set-DataRow($ws, $row, $runDate, $runTime, $sprintName, $total, $pass, $fail, $timeTaken) {
    Set-LiteralText $ws $row 1 $runDate
    Set-RunTime $ws $row $runTime
    Set-LiteralText $ws $row 3 $sprintName
    Set-PlainValue $ws $row 4 $total
    Set-PlainValue $ws $row 5 $pass
    Set-PlainValue $ws $row 6 $fail
    Set-PlainValue $ws $row 7 $timeTaken
}

# ---------------------------------------------------------------------
# AMSIN - previous last row (80) gains the standard row style, its
# Run Time gets corrected to the precise timestamp, then two new runs
# (81, 82) are appended.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AMSIN")

Set-LiteralText $ws 80 1 "2023-04-18"
Set-RunTime $ws 80 45034.57697421296
Set-LiteralText $ws 80 3 "176firsttrail"
Set-PlainValue $ws 80 4 89
Set-PlainValue $ws 80 5 89
Set-PlainValue $ws 80 6 0
Set-PlainValue $ws 80 7 2.57

Set-DataRow $ws 81 "2023-04-19" 45035.70495701389 "176scndcyc" 89 89 0 2.56
Set-DataRow $ws 82 "2023-04-20" 45036.41587881945 "176fnlruntest" 89 89 0 2.57

# ---------------------------------------------------------------------
# BETA - one new run (34) appended.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BETA")

Set-DataRow $ws 34 "2023-04-20" 45036.52489837963 "176beta" 89 89 0 2.43

# ---------------------------------------------------------------------
# AMS - two new runs (46, 47) appended.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AMS")

Set-DataRow $ws 46 "2023-05-04" 45050.71928640046 "176firstsycle" 89 89 0 2.3
Set-DataRow $ws 47 "2023-05-08" 45054.54537790914 "176htfxtrl" 89 89 0 2.17
